# Update the "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# Both sheets share identical data for rows 2-20; apply the same updates to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1037
    3  = 311
    4  = 1423
    5  = 8575
    6  = 70
    9  = 262
    10 = 148
    11 = 3456
    13 = 348
    14 = 71
    15 = 1019
    18 = 301
    19 = 175
    20 = 2152
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
